# Fruta / hortaliza, semanal
#
# The rows 3-14 (the weekly Haba / Mapocho Venta Directa records) get their
# per-record fields (Fecha, Volumen, Precio mínimo/máximo/promedio, Origen,
# Precio $/Kg) reshuffled across the date-ordered rows while the
# market/region/category/quality/unit/kg columns stay put. Capture a
# snapshot of the "before" values for the columns that move, then write
# them back out according to the row permutation implied by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) that get shuffled between rows 3..14:
#   D = 4 (Fecha), J = 10 (Volumen), K = 11 (Precio mínimo),
#   L = 12 (Precio máximo), M = 13 (Precio promedio ponderado),
#   O = 15 (Origen), P = 16 (Precio $/Kg)
$cols = @(4, 10, 11, 12, 13, 15, 16)

# Snapshot current ("before") values for rows 3..14 across those columns.
$snapshot = @{}
foreach ($r in 3..14) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $row
}

# Target row <- source row permutation (source row's old values are
# written into the target row).
$perm = @{
    3  = 10
    4  = 6
    5  = 11
    6  = 14
    7  = 12
    8  = 5
    9  = 4
    10 = 13
    11 = 3
    12 = 9
    13 = 7
    14 = 8
}

foreach ($targetRow in 3..14) {
    $sourceRow = $perm[$targetRow]
    $srcValues = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcValues[$c]
    }
}
